$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.461.29'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.878.85'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.04'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +5.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4770'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.96%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06523'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.91'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07726'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.22'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7382'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +8.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.877.18'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.132'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '273.63'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.443.57'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.61'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007584'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9999'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.124.20'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.258'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.186'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.337'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.98'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.86'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.946'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.372'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09956'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.519'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.312'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.068'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04792'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.127'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7011'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.714'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.727'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.338'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '71.30'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.952'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4212'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8362'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.95'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.273'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.086'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.64'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '931.55'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05645'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.31%  '
